$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 (fedot)
$ws.Range("B8").Value = "0.989 (0.988 Â± 0.000)"
$ws.Range("C8").Value = "00:05:00 (00:06:51 Â± 00:01:55)"
$ws.Range("E8").Value = "[2, 3, 5, 11, 13, 17, 19, 23, 29, 31, 37, 41, 47, 59, 61, 67, 71]"

# Row 16 (pycaret)
$ws.Range("B16").Value = "1.000 (0.995 Â± 0.005)"
$ws.Range("C16").Value = "00:00:09 (00:00:10 Â± 00:00:00)"
